$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct week 6 (row 7) case count
$ws.Range("B7").Value = 454

# Add week 7 (row 8) data
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 13
